$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "1.001") are pre-formatted
# as Text so Excel stores them as strings (matching the source
# inline-string cells) instead of coercing them to doubles.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '23.341.55'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.627.27'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("D6").Value = '303.57'
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("D7").Value = '0.3774'
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").Value = '51.96'
$ws.Range("E8").Value = '  -2.57%  '
$ws.Range("D9").Value = '0.3618'
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("D11").Value = '0.08077'
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '22.63'
$ws.Range("E13").Value = '  -2.47%  '
$ws.Range("D14").Value = '6.559'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '0.00001243'
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '7.222'
$ws.Range("E16").Value = '  -3.22%  '
$ws.Range("D17").Value = '1.627.64'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '93.45'
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").Value = '0.06912'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").Value = '17.93'
$ws.Range("E20").Value = '  -2.51%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").Value = '6.443'
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("D23").Value = '23.356.22'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("D25").Value = '3.222'
$ws.Range("E25").Value = '  +2.93%  '
$ws.Range("D26").Value = '2.455'
$ws.Range("E26").Value = '  +1.32%  '
$ws.Range("D27").Value = '21.05'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '148.86'
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").Value = '5.288'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Value = '134.73'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("D31").Value = '2.295'
$ws.Range("E31").Value = '  -4.85%  '
$ws.Range("D32").Value = '1.809.64'
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("D33").Value = '6.752'
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").Value = '10.92'
$ws.Range("E34").Value = '  +4.22%  '
$ws.Range("D35").Value = '0.9461'
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("D36").Value = '0.02815'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = '0.2528'
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").Value = '0.08812'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '6.108'
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("D40").Value = '0.07098'
$ws.Range("E40").Value = '  -4.86%  '
$ws.Range("D41").Value = '1.362'
$ws.Range("E41").Value = '  -3.35%  '
$ws.Range("D42").Value = '0.7042'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").Value = '12.33'
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("D45").Value = '0.6447'
$ws.Range("E45").Value = '  -2.60%  '
$ws.Range("D46").Value = '2.316'
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").Value = '3.979'
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").Value = '0.07975'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").Value = '1.203'
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").Value = '125.66'
$ws.Range("E51").Value = '  -4.31%  '

# Restore the default (unstyled) cell format now that the text values
# are committed, so styling matches the original workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
